# Update the "team specific time data" probability matrix cells
# (values recomputed from refreshed team-specific timing stats).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2071611253196931
$ws.Range("C2").Value = 0.5115089514066496
$ws.Range("J2").Value = 0.03324808184143223
$ws.Range("O2").Value = 0.002557544757033248
$ws.Range("P2").Value = 0.1483375959079284
$ws.Range("S2").Value = 0.09718670076726342

# Row 3
$ws.Range("B3").Value = 0.0182648401826484
$ws.Range("C3").Value = 0.0091324200913242
$ws.Range("J3").Value = 0.0273972602739726
$ws.Range("P3").Value = 0.7671232876712328
$ws.Range("S3").Value = 0.1780821917808219

# Row 4
$ws.Range("J4").Value = 0.1071428571428571
$ws.Range("P4").Value = 0.6785714285714286
$ws.Range("S4").Value = 0.2142857142857143

# Row 6
$ws.Range("B6").Value = 0.05882352941176471
$ws.Range("D6").Value = 0.02941176470588235
$ws.Range("F6").Value = 0.07563025210084033
$ws.Range("J6").Value = 0.2815126050420168
$ws.Range("O6").Value = 0.02100840336134454
$ws.Range("Q6").Value = 0.1260504201680672
$ws.Range("R6").Value = 0.09243697478991597
$ws.Range("S6").Value = 0.3151260504201681

# Row 7
$ws.Range("B7").Value = 0.1318181818181818
$ws.Range("D7").Value = 0.02272727272727273
$ws.Range("F7").Value = 0.06818181818181818
$ws.Range("J7").Value = 0.1181818181818182
$ws.Range("O7").Value = 0.01818181818181818
$ws.Range("Q7").Value = 0.1681818181818182
$ws.Range("R7").Value = 0.06818181818181818
$ws.Range("S7").Value = 0.4045454545454545

# Row 8
$ws.Range("B8").Value = 0.128
$ws.Range("D8").Value = 0.016
$ws.Range("F8").Value = 0.074
$ws.Range("J8").Value = 0.134
$ws.Range("O8").Value = 0.016
$ws.Range("Q8").Value = 0.152
$ws.Range("R8").Value = 0.114
$ws.Range("S8").Value = 0.366

# Row 9
$ws.Range("B9").Value = 0.09677419354838709
$ws.Range("D9").Value = 0.01935483870967742
$ws.Range("F9").Value = 0.06451612903225806
$ws.Range("J9").Value = 0.09677419354838709
$ws.Range("O9").Value = 0.01290322580645161
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.07741935483870968
$ws.Range("S9").Value = 0.432258064516129

# Row 10
$ws.Range("B10").Value = 0.1283185840707965
$ws.Range("D10").Value = 0.02359882005899705
$ws.Range("F10").Value = 0.05752212389380531
$ws.Range("J10").Value = 0.1541297935103245
$ws.Range("O10").Value = 0.01991150442477876
$ws.Range("Q10").Value = 0.1902654867256637
$ws.Range("R10").Value = 0.0803834808259587
$ws.Range("S10").Value = 0.3458702064896755

# Row 11
$ws.Range("G11").Value = 0.1436950146627566
$ws.Range("J11").Value = 0.1055718475073314
$ws.Range("K11").Value = 0.2170087976539589
$ws.Range("L11").Value = 0.5219941348973607
$ws.Range("S11").Value = 0.01173020527859238

# Row 12
$ws.Range("G12").Value = 0.7513812154696132
$ws.Range("J12").Value = 0.1878453038674033
$ws.Range("K12").Value = 0.01104972375690608
$ws.Range("L12").Value = 0.01104972375690608
$ws.Range("S12").Value = 0.03867403314917127

# Row 13
$ws.Range("F13").Value = 0.01515151515151515
$ws.Range("G13").Value = 0.6363636363636364
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.07575757575757576

# Row 15
$ws.Range("F15").Value = 0.02127659574468085
$ws.Range("H15").Value = 0.1617021276595745
$ws.Range("I15").Value = 0.07234042553191489
$ws.Range("J15").Value = 0.3404255319148936
$ws.Range("K15").Value = 0.1063829787234043
$ws.Range("M15").Value = 0.008510638297872341
$ws.Range("O15").Value = 0.07659574468085106
$ws.Range("S15").Value = 0.2127659574468085

# Row 16
$ws.Range("F16").Value = 0.0196078431372549
$ws.Range("H16").Value = 0.1529411764705882
$ws.Range("I16").Value = 0.0392156862745098
$ws.Range("J16").Value = 0.4392156862745098
$ws.Range("K16").Value = 0.1568627450980392
$ws.Range("M16").Value = 0.03137254901960784
$ws.Range("N16").Value = 0.00392156862745098
$ws.Range("O16").Value = 0.05490196078431372
$ws.Range("S16").Value = 0.1019607843137255

# Row 17
$ws.Range("F17").Value = 0.03167420814479638
$ws.Range("H17").Value = 0.1764705882352941
$ws.Range("I17").Value = 0.083710407239819
$ws.Range("J17").Value = 0.4208144796380091
$ws.Range("K17").Value = 0.1063348416289593
$ws.Range("M17").Value = 0.0248868778280543
$ws.Range("N17").Value = 0.002262443438914027
$ws.Range("O17").Value = 0.05429864253393665
$ws.Range("S17").Value = 0.09954751131221719

# Row 18
$ws.Range("F18").Value = 0.0182648401826484
$ws.Range("H18").Value = 0.228310502283105
$ws.Range("I18").Value = 0.0684931506849315
$ws.Range("J18").Value = 0.3698630136986301
$ws.Range("K18").Value = 0.091324200913242
$ws.Range("M18").Value = 0.0273972602739726
$ws.Range("O18").Value = 0.0776255707762557
$ws.Range("S18").Value = 0.1187214611872146

# Row 19
$ws.Range("F19").Value = 0.01998401278976818
$ws.Range("H19").Value = 0.2326139088729017
$ws.Range("I19").Value = 0.05755395683453238
$ws.Range("J19").Value = 0.3565147881694644
$ws.Range("K19").Value = 0.1135091926458833
$ws.Range("M19").Value = 0.03037569944044764
$ws.Range("N19").Value = 0.0007993605115907274
$ws.Range("O19").Value = 0.07434052757793765
$ws.Range("S19").Value = 0.114308553157474

Write-Host "Updated 111 cells in the team-specific time matrix."